$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.057.61"
$ws.Range("E2").Value = "  -2.18%  "

$ws.Range("D3").Value = "2.506.60"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'570.89"
$ws.Range("E5").Value = "  -0.88%  "

$ws.Range("D6").Value = "'166.13"
$ws.Range("E6").Value = "  -2.21%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D9").Value = "2.505.90"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("D12").Value = "'0.352"
$ws.Range("E12").Value = "  +2.70%  "

$ws.Range("D13").Value = "'4.91"
$ws.Range("E13").Value = "  +2.20%  "

$ws.Range("D14").Value = "2.976.94"
$ws.Range("E14").Value = "  -0.30%  "

$ws.Range("D15").Value = "68.979.51"
$ws.Range("E15").Value = "  -2.20%  "

$ws.Range("E16").Value = "  -3.12%  "

$ws.Range("D17").Value = "'24.79"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "2.503.24"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("D19").Value = "'11.31"
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").Value = "'348.25"
$ws.Range("E21").Value = "  -2.04%  "

$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").Value = "'1.97"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'70.10"
$ws.Range("E25").Value = "  +1.32%  "

$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("E27").Value = "  -2.86%  "

$ws.Range("D28").Value = "2.645.06"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "0.0₃0886"
$ws.Range("E30").Value = "  -2.58%  "

$ws.Range("D31").Value = "'7.83"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").Value = "'460.43"
$ws.Range("E32").Value = "  -3.70%  "

$ws.Range("D33").Value = "'1.23"
$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  +1.38%  "

$ws.Range("D37").Value = "'157.80"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("E43").Value = "  -3.19%  "

$ws.Range("D44").Value = "'38.23"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D46").Value = "'2.25"
$ws.Range("E46").Value = "  -6.39%  "

$ws.Range("D47").Value = "'141.25"
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").Value = "'0.525"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").Value = "'3.47"
$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("D50").Value = "'0.0728"
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").Value = "'1.55"
$ws.Range("E51").Value = "  -3.53%  "
